$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H2").Value = 349.66666
$ws_ALC.Range("I2").Value = 325
$ws_ALC.Range("J2").Value = 399
$ws_ALC.Range("K2").Value = 325
$ws_ALC.Range("L2").Value = 399
$ws_ALC.Range("M2").Value = -212
$ws_ALC.Range("N2").Value = -625
$ws_ALC.Range("H4").Value = 70.5
$ws_ALC.Range("I4").Value = 99
$ws_ALC.Range("J4").Value = 42
$ws_ALC.Range("K4").Value = 99
$ws_ALC.Range("L4").Value = 42
$ws_ALC.Range("M4").Value = 15
$ws_ALC.Range("N4").Value = -270
$ws_ALC.Range("H6").Value = 334.2143
$ws_ALC.Range("I6").Value = 334.2143
$ws_ALC.Range("J6").Value = 0
$ws_ALC.Range("K6").Value = 1002.6429
$ws_ALC.Range("L6").Value = 0
$ws_ALC.Range("M6").Value = -890.6428999999999
$ws_ALC.Range("N6").ClearContents()
$ws_ALC.Range("H43").Value = 697.3333
$ws_ALC.Range("I43").Value = 0
$ws_ALC.Range("J43").Value = 697.3333
$ws_ALC.Range("K43").Value = 0
$ws_ALC.Range("L43").Value = 697.3333
$ws_ALC.Range("N43").Value = -835.3333
$ws_ALC.Range("H64").Value = 3980.2068
$ws_ALC.Range("I64").Value = 3093.5293
$ws_ALC.Range("J64").Value = 5236.3335
$ws_ALC.Range("K64").Value = 3093.5293
$ws_ALC.Range("L64").Value = 5236.3335
$ws_ALC.Range("M64").Value = -2845.5293
$ws_ALC.Range("N64").Value = -5732.3335
$ws_ALC.Range("H67").Value = 3980.2068
$ws_ALC.Range("I67").Value = 3093.5293
$ws_ALC.Range("J67").Value = 5236.3335
$ws_ALC.Range("K67").Value = 3093.5293
$ws_ALC.Range("L67").Value = 5236.3335
$ws_ALC.Range("M67").Value = -2235.5293
$ws_ALC.Range("N67").Value = -6952.3335
$ws_ALC.Range("H76").Value = 7939993.5
$ws_ALC.Range("I76").Value = 15876087
$ws_ALC.Range("J76").Value = 3900
$ws_ALC.Range("K76").Value = 15876087
$ws_ALC.Range("L76").Value = 3900
$ws_ALC.Range("M76").Value = -15875772
$ws_ALC.Range("N76").Value = -4530
$ws_ALC.Range("H79").Value = 7939993.5
$ws_ALC.Range("I79").Value = 15876087
$ws_ALC.Range("J79").Value = 3900
$ws_ALC.Range("K79").Value = 15876087
$ws_ALC.Range("L79").Value = 3900
$ws_ALC.Range("M79").Value = -15874995
$ws_ALC.Range("N79").Value = -6084
$ws_ALC.Range("H88").Value = 13892700
$ws_ALC.Range("I88").Value = 2000
$ws_ALC.Range("J88").Value = 15877085
$ws_ALC.Range("K88").Value = 2000
$ws_ALC.Range("L88").Value = 15877085
$ws_ALC.Range("M88").Value = -1594
$ws_ALC.Range("N88").Value = -15877897
$ws_ALC.Range("H91").Value = 13892700
$ws_ALC.Range("I91").Value = 2000
$ws_ALC.Range("J91").Value = 15877085
$ws_ALC.Range("K91").Value = 2000
$ws_ALC.Range("L91").Value = 15877085
$ws_ALC.Range("M91").Value = -596
$ws_ALC.Range("N91").Value = -15879893
$ws_ALC.Range("H92").Value = 556061.9
$ws_ALC.Range("I92").Value = 585307.2
$ws_ALC.Range("J92").Value = 400
$ws_ALC.Range("K92").Value = 585307.2
$ws_ALC.Range("L92").Value = 400
$ws_ALC.Range("M92").Value = -584059.2
$ws_ALC.Range("N92").Value = -2896
$ws_ALC.Range("H94").Value = 11409.546
$ws_ALC.Range("I94").Value = 12389.444
$ws_ALC.Range("J94").Value = 7000
$ws_ALC.Range("K94").Value = 12389.444
$ws_ALC.Range("L94").Value = 7000
$ws_ALC.Range("M94").Value = -11938.444
$ws_ALC.Range("N94").Value = -7902
$ws_ALC.Range("H98").Value = 745381.25
$ws_ALC.Range("I98").Value = 859132.5600000001
$ws_ALC.Range("J98").Value = 5998
$ws_ALC.Range("K98").Value = 859132.5600000001
$ws_ALC.Range("L98").Value = 5998
$ws_ALC.Range("M98").Value = -857634.5600000001
$ws_ALC.Range("N98").Value = -8994
$ws_ALC.Range("H101").Value = 455.14285
$ws_ALC.Range("I101").Value = 381.0909
$ws_ALC.Range("J101").Value = 726.6667
$ws_ALC.Range("K101").Value = 1143.2727
$ws_ALC.Range("L101").Value = 2180.0001
$ws_ALC.Range("M101").Value = 478.7273
$ws_ALC.Range("N101").Value = -5424.0001
$ws_ALC.Range("H106").Value = 27778854
$ws_ALC.Range("I106").Value = 27778854
$ws_ALC.Range("J106").Value = 0
$ws_ALC.Range("K106").Value = 27778854
$ws_ALC.Range("L106").Value = 0
$ws_ALC.Range("M106").Value = -27778223
$ws_ALC.Range("H122").Value = 745381.25
$ws_ALC.Range("I122").Value = 859132.5600000001
$ws_ALC.Range("J122").Value = 5998
$ws_ALC.Range("K122").Value = 2577397.68
$ws_ALC.Range("L122").Value = 17994
$ws_ALC.Range("M122").Value = -2574947.68
$ws_ALC.Range("N122").Value = -22894

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H6").Value = 2537501.5
$ws_ARM.Range("I6").Value = 2537501.5
$ws_ARM.Range("J6").Value = 0
$ws_ARM.Range("K6").Value = 2537501.5
$ws_ARM.Range("L6").Value = 0
$ws_ARM.Range("M6").Value = -2537328.5
$ws_ARM.Range("N6").ClearContents()
$ws_ARM.Range("H26").Value = 0
$ws_ARM.Range("I26").Value = 0
$ws_ARM.Range("J26").Value = 0
$ws_ARM.Range("K26").Value = 0
$ws_ARM.Range("L26").Value = 0
$ws_ARM.Range("M26").ClearContents()
$ws_ARM.Range("N26").ClearContents()
$ws_ARM.Range("H61").Value = 3405.24
$ws_ARM.Range("I61").Value = 2528.4
$ws_ARM.Range("J61").Value = 4720.5
$ws_ARM.Range("K61").Value = 2528.4
$ws_ARM.Range("L61").Value = 4720.5
$ws_ARM.Range("M61").Value = -2316.4
$ws_ARM.Range("N61").Value = -5144.5
$ws_ARM.Range("H132").Value = 2465.8276
$ws_ARM.Range("I132").Value = 2043.8776
$ws_ARM.Range("J132").Value = 4763.1113
$ws_ARM.Range("K132").Value = 6131.6328
$ws_ARM.Range("L132").Value = 14289.3339
$ws_ARM.Range("M132").Value = -3601.6328
$ws_ARM.Range("N132").Value = -19349.3339
$ws_ARM.Range("H136").Value = 3405.24
$ws_ARM.Range("I136").Value = 2528.4
$ws_ARM.Range("J136").Value = 4720.5
$ws_ARM.Range("K136").Value = 7585.200000000001
$ws_ARM.Range("L136").Value = 14161.5
$ws_ARM.Range("M136").Value = -5035.200000000001
$ws_ARM.Range("N136").Value = -19261.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H82").Value = 8000
$ws_BSM.Range("I82").Value = 8000
$ws_BSM.Range("J82").Value = 0
$ws_BSM.Range("K82").Value = 8000
$ws_BSM.Range("L82").Value = 0
$ws_BSM.Range("M82").Value = -7617
$ws_BSM.Range("N82").ClearContents()
$ws_BSM.Range("H85").Value = 8000
$ws_BSM.Range("I85").Value = 8000
$ws_BSM.Range("J85").Value = 0
$ws_BSM.Range("K85").Value = 8000
$ws_BSM.Range("L85").Value = 0
$ws_BSM.Range("M85").Value = -6674
$ws_BSM.Range("N85").ClearContents()
$ws_BSM.Range("H86").Value = 6960.55
$ws_BSM.Range("I86").Value = 3030.5
$ws_BSM.Range("J86").Value = 8644.857
$ws_BSM.Range("K86").Value = 3030.5
$ws_BSM.Range("L86").Value = 8644.857
$ws_BSM.Range("M86").Value = -1907.5
$ws_BSM.Range("N86").Value = -10890.857
$ws_BSM.Range("H89").Value = 6960.55
$ws_BSM.Range("I89").Value = 3030.5
$ws_BSM.Range("J89").Value = 8644.857
$ws_BSM.Range("K89").Value = 15152.5
$ws_BSM.Range("L89").Value = 43224.285
$ws_BSM.Range("M89").Value = -9536.5
$ws_BSM.Range("N89").Value = -54456.285
$ws_BSM.Range("H94").Value = 1438.6111
$ws_BSM.Range("I94").Value = 961
$ws_BSM.Range("J94").Value = 3826.6667
$ws_BSM.Range("K94").Value = 961
$ws_BSM.Range("L94").Value = 3826.6667
$ws_BSM.Range("M94").Value = -510
$ws_BSM.Range("N94").Value = -4728.6667

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H12").Value = 10602
$ws_CRP.Range("I12").Value = 13127.5
$ws_CRP.Range("J12").Value = 500
$ws_CRP.Range("K12").Value = 13127.5
$ws_CRP.Range("L12").Value = 500
$ws_CRP.Range("M12").Value = -12957.5
$ws_CRP.Range("N12").Value = -840
$ws_CRP.Range("H58").Value = 2532.9048
$ws_CRP.Range("I58").Value = 1740.4286
$ws_CRP.Range("J58").Value = 4117.857
$ws_CRP.Range("K58").Value = 1740.4286
$ws_CRP.Range("L58").Value = 4117.857
$ws_CRP.Range("M58").Value = -1537.4286
$ws_CRP.Range("N58").Value = -4523.857
$ws_CRP.Range("H132").Value = 2786.4
$ws_CRP.Range("I132").Value = 2342.182
$ws_CRP.Range("J132").Value = 4880.5713
$ws_CRP.Range("K132").Value = 7026.545999999999
$ws_CRP.Range("L132").Value = 14641.7139
$ws_CRP.Range("M132").Value = -4496.545999999999
$ws_CRP.Range("N132").Value = -19701.7139
$ws_CRP.Range("H136").Value = 2532.9048
$ws_CRP.Range("I136").Value = 1740.4286
$ws_CRP.Range("J136").Value = 4117.857
$ws_CRP.Range("K136").Value = 5221.2858
$ws_CRP.Range("L136").Value = 12353.571
$ws_CRP.Range("M136").Value = -2671.2858
$ws_CRP.Range("N136").Value = -17453.571

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H6").Value = 436.91666
$ws_CUL.Range("I6").Value = 77.28570999999999
$ws_CUL.Range("J6").Value = 940.4
$ws_CUL.Range("K6").Value = 231.85713
$ws_CUL.Range("L6").Value = 2821.2
$ws_CUL.Range("M6").Value = -118.85713
$ws_CUL.Range("N6").Value = -3047.2
$ws_CUL.Range("H75").Value = 1641
$ws_CUL.Range("I75").Value = 784.75
$ws_CUL.Range("J75").Value = 1904.4615
$ws_CUL.Range("K75").Value = 2354.25
$ws_CUL.Range("L75").Value = 5713.3845
$ws_CUL.Range("M75").Value = -1356.25
$ws_CUL.Range("N75").Value = -7709.3845
$ws_CUL.Range("H78").Value = 1641
$ws_CUL.Range("I78").Value = 784.75
$ws_CUL.Range("J78").Value = 1904.4615
$ws_CUL.Range("K78").Value = 7062.75
$ws_CUL.Range("L78").Value = 17140.1535
$ws_CUL.Range("M78").Value = -2070.75
$ws_CUL.Range("N78").Value = -27124.1535

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H15").Value = 0
$ws_GSM.Range("I15").Value = 0
$ws_GSM.Range("J15").Value = 0
$ws_GSM.Range("K15").Value = 0
$ws_GSM.Range("L15").Value = 0
$ws_GSM.Range("N15").ClearContents()
$ws_GSM.Range("H81").Value = 0
$ws_GSM.Range("I81").Value = 0
$ws_GSM.Range("J81").Value = 0
$ws_GSM.Range("K81").Value = 0
$ws_GSM.Range("L81").Value = 0
$ws_GSM.Range("N81").ClearContents()
$ws_GSM.Range("H84").Value = 0
$ws_GSM.Range("I84").Value = 0
$ws_GSM.Range("J84").Value = 0
$ws_GSM.Range("K84").Value = 0
$ws_GSM.Range("L84").Value = 0
$ws_GSM.Range("N84").ClearContents()
$ws_GSM.Range("H102").Value = 2449.3845
$ws_GSM.Range("I102").Value = 2635.625
$ws_GSM.Range("J102").Value = 2151.4
$ws_GSM.Range("K102").Value = 2635.625
$ws_GSM.Range("L102").Value = 2151.4
$ws_GSM.Range("M102").Value = -1013.625
$ws_GSM.Range("N102").Value = -5395.4
$ws_GSM.Range("H122").Value = 507132.53
$ws_GSM.Range("I122").Value = 856353.0600000001
$ws_GSM.Range("J122").Value = 2702.889
$ws_GSM.Range("K122").Value = 2569059.18
$ws_GSM.Range("L122").Value = 8108.667
$ws_GSM.Range("M122").Value = -2566609.18
$ws_GSM.Range("N122").Value = -13008.667
$ws_GSM.Range("H132").Value = 3872.8333
$ws_GSM.Range("I132").Value = 3602.3
$ws_GSM.Range("J132").Value = 4211
$ws_GSM.Range("K132").Value = 10806.9
$ws_GSM.Range("L132").Value = 12633
$ws_GSM.Range("M132").Value = -8276.900000000001
$ws_GSM.Range("N132").Value = -17693

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 3078.7666
$ws_LTW.Range("I40").Value = 2712.6155
$ws_LTW.Range("J40").Value = 3358.7646
$ws_LTW.Range("K40").Value = 2712.6155
$ws_LTW.Range("L40").Value = 3358.7646
$ws_LTW.Range("M40").Value = -2576.6155
$ws_LTW.Range("N40").Value = -3630.7646
$ws_LTW.Range("H46").Value = 1193.5
$ws_LTW.Range("I46").Value = 761.25
$ws_LTW.Range("J46").Value = 1539.3
$ws_LTW.Range("K46").Value = 761.25
$ws_LTW.Range("L46").Value = 1539.3
$ws_LTW.Range("M46").Value = -573.25
$ws_LTW.Range("N46").Value = -1915.3
$ws_LTW.Range("H101").Value = 22154.285
$ws_LTW.Range("I101").Value = 0
$ws_LTW.Range("J101").Value = 22154.285
$ws_LTW.Range("K101").Value = 0
$ws_LTW.Range("L101").Value = 22154.285
$ws_LTW.Range("N101").Value = -28644.285
$ws_LTW.Range("H122").Value = 2893.5356
$ws_LTW.Range("I122").Value = 1967.9333
$ws_LTW.Range("J122").Value = 3961.5386
$ws_LTW.Range("K122").Value = 5903.7999
$ws_LTW.Range("L122").Value = 11884.6158
$ws_LTW.Range("M122").Value = -3453.7999
$ws_LTW.Range("N122").Value = -16784.6158
$ws_LTW.Range("H133").Value = 46759.273
$ws_LTW.Range("I133").Value = 0
$ws_LTW.Range("J133").Value = 46759.273
$ws_LTW.Range("K133").Value = 0
$ws_LTW.Range("L133").Value = 46759.273
$ws_LTW.Range("N133").Value = -51819.273
$ws_LTW.Range("H137").Value = 0
$ws_LTW.Range("I137").Value = 0
$ws_LTW.Range("J137").Value = 0
$ws_LTW.Range("K137").Value = 0
$ws_LTW.Range("L137").Value = 0
$ws_LTW.Range("N137").ClearContents()
